$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "40+12=52"
$t.Cell(1,2).Range.Text = "88-80=8"
$t.Cell(1,3).Range.Text = "42+56=98"
$t.Cell(1,4).Range.Text = "48-42=6"
$t.Cell(1,5).Range.Text = "51-43=8"
$t.Cell(2,1).Range.Text = "5+13=18"
$t.Cell(2,2).Range.Text = "66-58=8"
$t.Cell(2,3).Range.Text = "36+31=67"
$t.Cell(2,4).Range.Text = "46-22=24"
$t.Cell(2,5).Range.Text = "48+36=84"
$t.Cell(3,1).Range.Text = "68+7=75"
$t.Cell(3,2).Range.Text = "68-28=40"
$t.Cell(3,3).Range.Text = "18+79=97"
$t.Cell(3,4).Range.Text = "92-10=82"
$t.Cell(3,5).Range.Text = "86-68=18"
$t.Cell(4,1).Range.Text = "84-27=57"
$t.Cell(4,2).Range.Text = "29+58=87"
$t.Cell(4,3).Range.Text = "1+46=47"
$t.Cell(4,4).Range.Text = "80-7=73"
$t.Cell(4,5).Range.Text = "11-11=0"
$t.Cell(5,1).Range.Text = "67+19=86"
$t.Cell(5,2).Range.Text = "6+30=36"
$t.Cell(5,3).Range.Text = "61-52=9"
$t.Cell(5,4).Range.Text = "37-3=34"
$t.Cell(5,5).Range.Text = "20-1=19"
$t.Cell(6,1).Range.Text = "1+91=92"
$t.Cell(6,2).Range.Text = "16+18=34"
$t.Cell(6,3).Range.Text = "41-7=34"
$t.Cell(6,4).Range.Text = "60+17=77"
$t.Cell(6,5).Range.Text = "21+40=61"
$t.Cell(7,1).Range.Text = "76-47=29"
$t.Cell(7,2).Range.Text = "54-1=53"
$t.Cell(7,3).Range.Text = "97-38=59"
$t.Cell(7,4).Range.Text = "12-0=12"
$t.Cell(7,5).Range.Text = "24+61=85"
$t.Cell(8,1).Range.Text = "7+69=76"
$t.Cell(8,2).Range.Text = "92-18=74"
$t.Cell(8,3).Range.Text = "51-21=30"
$t.Cell(8,4).Range.Text = "62-37=25"
$t.Cell(8,5).Range.Text = "72-59=13"
$t.Cell(9,1).Range.Text = "48-0=48"
$t.Cell(9,2).Range.Text = "19+30=49"
$t.Cell(9,3).Range.Text = "46-34=12"
$t.Cell(9,4).Range.Text = "31+35=66"
$t.Cell(9,5).Range.Text = "37-6=31"
$t.Cell(10,1).Range.Text = "88-27=61"
$t.Cell(10,2).Range.Text = "53+18=71"
$t.Cell(10,3).Range.Text = "44+41=85"
$t.Cell(10,4).Range.Text = "40+11=51"
$t.Cell(10,5).Range.Text = "46+25=71"
$t.Cell(11,1).Range.Text = "1+60=61"
$t.Cell(11,2).Range.Text = "59+40=99"
$t.Cell(11,3).Range.Text = "21+2=23"
$t.Cell(11,4).Range.Text = "33-27=6"
$t.Cell(11,5).Range.Text = "70+15=85"
$t.Cell(12,1).Range.Text = "2+67=69"
$t.Cell(12,2).Range.Text = "30+15=45"
$t.Cell(12,3).Range.Text = "76-61=15"
$t.Cell(12,4).Range.Text = "92+1=93"
$t.Cell(12,5).Range.Text = "71+18=89"
$t.Cell(13,1).Range.Text = "93-38=55"
$t.Cell(13,2).Range.Text = "34+8=42"
$t.Cell(13,3).Range.Text = "24-0=24"
$t.Cell(13,4).Range.Text = "80-25=55"
$t.Cell(13,5).Range.Text = "36+54=90"
$t.Cell(14,1).Range.Text = "38-31=7"
$t.Cell(14,2).Range.Text = "2+69=71"
$t.Cell(14,3).Range.Text = "82-76=6"
$t.Cell(14,4).Range.Text = "9+89=98"
$t.Cell(14,5).Range.Text = "67-1=66"
$t.Cell(15,1).Range.Text = "31+59=90"
$t.Cell(15,2).Range.Text = "77-56=21"
$t.Cell(15,3).Range.Text = "8+57=65"
$t.Cell(15,4).Range.Text = "26+6=32"
$t.Cell(15,5).Range.Text = "16+82=98"
$t.Cell(16,1).Range.Text = "81-49=32"
$t.Cell(16,2).Range.Text = "51+41=92"
$t.Cell(16,3).Range.Text = "21+19=40"
$t.Cell(16,4).Range.Text = "70+25=95"
$t.Cell(16,5).Range.Text = "21+38=59"
$t.Cell(17,1).Range.Text = "50+46=96"
$t.Cell(17,2).Range.Text = "24+30=54"
$t.Cell(17,3).Range.Text = "29-1=28"
$t.Cell(17,4).Range.Text = "40-24=16"
$t.Cell(17,5).Range.Text = "8+0=8"
$t.Cell(18,1).Range.Text = "94-56=38"
$t.Cell(18,2).Range.Text = "10-2=8"
$t.Cell(18,3).Range.Text = "25+43=68"
$t.Cell(18,4).Range.Text = "83-27=56"
$t.Cell(18,5).Range.Text = "19-11=8"
$t.Cell(19,1).Range.Text = "26+67=93"
$t.Cell(19,2).Range.Text = "95-22=73"
$t.Cell(19,3).Range.Text = "18+12=30"
$t.Cell(19,4).Range.Text = "59+5=64"
$t.Cell(19,5).Range.Text = "4+0=4"
$t.Cell(20,1).Range.Text = "94+3=97"
$t.Cell(20,2).Range.Text = "75-52=23"
$t.Cell(20,3).Range.Text = "1+12=13"
$t.Cell(20,4).Range.Text = "81-26=55"
$t.Cell(20,5).Range.Text = "70-10=60"
